# Append the next day (serials 45761.x, half-hourly) of RES power trading price
# data (date, value, price) to Sheet1, continuing directly after the existing
# rows 2-49. This is step 1 of wiring up the automatic GMAIL-API price export.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(50, 1).Value = 45761
$ws.Cells.Item(50, 2).Value = 23433.8
$ws.Cells.Item(50, 3).Value = 39.229999999999997

$ws.Cells.Item(51, 1).Value = 45761.041666666657
$ws.Cells.Item(51, 2).Value = 23721.599999999999
$ws.Cells.Item(51, 3).Value = 38.159999999999997

$ws.Cells.Item(52, 1).Value = 45761.083333333343
$ws.Cells.Item(52, 2).Value = 24885.4
$ws.Cells.Item(52, 3).Value = 37.74

$ws.Cells.Item(53, 1).Value = 45761.125
$ws.Cells.Item(53, 2).Value = 26720.400000000001
$ws.Cells.Item(53, 3).Value = 35

$ws.Cells.Item(54, 1).Value = 45761.166666666657
$ws.Cells.Item(54, 2).Value = 27088.2
$ws.Cells.Item(54, 3).Value = 30.01

$ws.Cells.Item(55, 1).Value = 45761.208333333343
$ws.Cells.Item(55, 2).Value = 25703.599999999999
$ws.Cells.Item(55, 3).Value = 35.01

$ws.Cells.Item(56, 1).Value = 45761.25
$ws.Cells.Item(56, 2).Value = 25638.6
$ws.Cells.Item(56, 3).Value = 42.5

$ws.Cells.Item(57, 1).Value = 45761.291666666657
$ws.Cells.Item(57, 2).Value = 25631
$ws.Cells.Item(57, 3).Value = 68.19

$ws.Cells.Item(58, 1).Value = 45761.333333333343
$ws.Cells.Item(58, 2).Value = 26469.4
$ws.Cells.Item(58, 3).Value = 74.94

$ws.Cells.Item(59, 1).Value = 45761.375
$ws.Cells.Item(59, 2).Value = 25453.200000000001
$ws.Cells.Item(59, 3).Value = 37.74

$ws.Cells.Item(60, 1).Value = 45761.416666666657
$ws.Cells.Item(60, 2).Value = 29166.400000000001
$ws.Cells.Item(60, 3).Value = 41.7

$ws.Cells.Item(61, 1).Value = 45761.458333333343
$ws.Cells.Item(61, 2).Value = 32201.599999999999
$ws.Cells.Item(61, 3).Value = 26.77

$ws.Cells.Item(62, 1).Value = 45761.5
$ws.Cells.Item(62, 2).Value = 34122
$ws.Cells.Item(62, 3).Value = 52.65

$ws.Cells.Item(63, 1).Value = 45761.541666666657
$ws.Cells.Item(63, 2).Value = 34824.400000000001
$ws.Cells.Item(63, 3).Value = 23.51

$ws.Cells.Item(64, 1).Value = 45761.583333333343
$ws.Cells.Item(64, 2).Value = 34221.800000000003
$ws.Cells.Item(64, 3).Value = 10.16

$ws.Cells.Item(65, 1).Value = 45761.625
$ws.Cells.Item(65, 2).Value = 30770.799999999999
$ws.Cells.Item(65, 3).Value = 20.239999999999998

$ws.Cells.Item(66, 1).Value = 45761.666666666657
$ws.Cells.Item(66, 2).Value = 27411.4
$ws.Cells.Item(66, 3).Value = 28.12

$ws.Cells.Item(67, 1).Value = 45761.708333333343
$ws.Cells.Item(67, 2).Value = 24996.2
$ws.Cells.Item(67, 3).Value = 38.46

$ws.Cells.Item(68, 1).Value = 45761.75
$ws.Cells.Item(68, 2).Value = 24532.799999999999
$ws.Cells.Item(68, 3).Value = 36.630000000000003

$ws.Cells.Item(69, 1).Value = 45761.791666666657
$ws.Cells.Item(69, 2).Value = 24587.200000000001
$ws.Cells.Item(69, 3).Value = 47.02

$ws.Cells.Item(70, 1).Value = 45761.833333333343
$ws.Cells.Item(70, 2).Value = 23119.200000000001
$ws.Cells.Item(70, 3).Value = 37.81

$ws.Cells.Item(71, 1).Value = 45761.875
$ws.Cells.Item(71, 2).Value = 22953.200000000001
$ws.Cells.Item(71, 3).Value = 35

$ws.Cells.Item(72, 1).Value = 45761.916666666657
$ws.Cells.Item(72, 2).Value = 20995
$ws.Cells.Item(72, 3).Value = 19.75

$ws.Cells.Item(73, 1).Value = 45761.958333333343
$ws.Cells.Item(73, 2).Value = 19951.8
$ws.Cells.Item(73, 3).Value = 13

# Column A keeps the "yyyy-mm-dd hh:mm:ss" date/time number format used by the rest
# of the date column (this reuses the existing style already applied to A2:A49).
$ws.Range("A50:A73").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# Scroll the sheet view down to the newly appended data and leave the selection
# where the user left off while reviewing it.
$ws.Range("M69").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
